$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit per-cell updates matching the authored diff.
# Numeric-looking strings (prices / percentages) are written with a
# leading apostrophe to force text, then re-styled "Normal" so the
# quote-prefix flag does not linger on the cell (matches the source
# workbook, which stores these as plain inline strings, no special style).

$ws.Range("D2").Value = "'283.33"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "'1.94%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'28.42"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "'4.25%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.054"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = "'3.68%"
$ws.Range("E4").Style = "Normal"

$ws.Range("E5").Value = "'0.95%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'7.226"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "'3.11%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'1.394"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = "'17.10%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9177"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "'3.79%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1536"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "'-1.51%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.06500"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "'27.24%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07613"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = "'1.64%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.02794"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = "'-3.20%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.08968"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = "'-0.05%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001586"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "'0.74%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006365"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = "'-0.56%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006063"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "'-1.06%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.444"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "'-1.07%"
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'1.47%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'-1.41%"
$ws.Range("E19").Style = "Normal"

$ws.Range("E21").Value = "'-0.59%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.008"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "'2.44%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.1544"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "'2.93%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.04429"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "'0.13%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'0.89%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.004461"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = "'15.07%"
$ws.Range("E26").Style = "Normal"

$ws.Range("E28").Value = "'1.65%"
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'-1.94%"
$ws.Range("E29").Style = "Normal"

$ws.Range("D40").Value = "'0.04116"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = "'-0.78%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.006614"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = "'-2.86%"
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'4.89%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002149"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "'13.70%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01152"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = "'2.04%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005394"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "'1.32%"
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "BOLO"

$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"

$ws.Range("D46").Value = "'1.933"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = "'14.55%"
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "CoinbaseStockToken"

$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"

$ws.Range("D47").Value = "'0.01851"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = "'-0.09%"
$ws.Range("E47").Style = "Normal"
